$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear header row (A1:R1) formatting: drop the bold/centered/bordered style ---
$ws.Range("A1:R1").Style = "Normal"

# --- Clear the "Unnamed: 0" label text from A1, leaving the cell blank (text-typed) ---
$ws.Range("A1").Value = "'"
$ws.Range("A1").Style = "Normal"

# --- Update numeric values across rows 3-8 (corrected pre/post/total fixation data) ---

# Row 3: Revisit count
$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 13
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 14
$ws.Range("I3").Value = 22
$ws.Range("J3").Value = 10
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 2
$ws.Range("O3").Value = "'"
$ws.Range("O3").Style = "Normal"
$ws.Range("P3").Value = 30
$ws.Range("Q3").Value = 4

# Row 4: Fixation count
$ws.Range("B4").Value = 69
$ws.Range("C4").Value = 25
$ws.Range("D4").Value = 20
$ws.Range("E4").Value = 27
$ws.Range("I4").Value = 53
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 4
$ws.Range("O4").Value = "'"
$ws.Range("O4").Style = "Normal"
$ws.Range("P4").Value = 159
$ws.Range("Q4").Value = 9

# Row 5: Dwell time (ms)
$ws.Range("B5").Value = 29533.32
$ws.Range("C5").Value = 13622.03
$ws.Range("D5").Value = 10268.63
$ws.Range("E5").Value = 14731.47
$ws.Range("I5").Value = 25751.21
$ws.Range("J5").Value = 8142.15
$ws.Range("K5").Value = 2077.05
$ws.Range("L5").Value = 1751.74
$ws.Range("O5").Value = "'"
$ws.Range("O5").Style = "Normal"
$ws.Range("P5").Value = 59677.22
$ws.Range("Q5").Value = 2986.89

# Row 6: Dwell time (%)
$ws.Range("B6").Value = 23.3
$ws.Range("C6").Value = 10.75
$ws.Range("D6").Value = 8.1
$ws.Range("E6").Value = 11.62
$ws.Range("I6").Value = 20.31
$ws.Range("J6").Value = 6.42
$ws.Range("K6").Value = 1.64
$ws.Range("L6").Value = 1.38
$ws.Range("M6").Value = 0.22
$ws.Range("O6").Value = "'"
$ws.Range("O6").Style = "Normal"
$ws.Range("P6").Value = 47.08
$ws.Range("Q6").Value = 2.36
$ws.Range("R6").Value = 0.22

# Row 7: Fixation duration (ms)
$ws.Range("B7").Value = 428.02
$ws.Range("C7").Value = 544.88
$ws.Range("D7").Value = 513.4299999999999
$ws.Range("E7").Value = 545.61
$ws.Range("I7").Value = 485.87
$ws.Range("J7").Value = 407.11
$ws.Range("K7").Value = 415.41
$ws.Range("L7").Value = 437.94
$ws.Range("O7").Value = "'"
$ws.Range("O7").Style = "Normal"
$ws.Range("P7").Value = 375.33
$ws.Range("Q7").Value = 331.88

# Row 8: First fixation duration (ms)
$ws.Range("O8").Value = "'"
$ws.Range("O8").Style = "Normal"

# --- Remove the trailing blank row 10 entirely (shrinks used range to A1:R9) ---
$ws.Rows("10").Delete()
